$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current last data row (row 74), pushing the
# existing row 74 down to row 75, then populate the new row 74 with the
# latest weekly record.
$ws.Rows.Item(74).Insert()

$ws.Cells.Item(74, 1).Value = 9
$ws.Cells.Item(74, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(74, 3).Value = "Metropolitana"
$ws.Cells.Item(74, 4).Value = 45121
$ws.Cells.Item(74, 5).Value = 13
$ws.Cells.Item(74, 6).Value = 100112010
$ws.Cells.Item(74, 7).Value = "Achicoria"
$ws.Cells.Item(74, 8).Value = "Sin especificar"
$ws.Cells.Item(74, 9).Value = "Primera"
$ws.Cells.Item(74, 10).Value = 70
$ws.Cells.Item(74, 11).Value = 8000
$ws.Cells.Item(74, 12).Value = 8000
$ws.Cells.Item(74, 13).Value = 8000
$ws.Cells.Item(74, 14).Value = "$/caja 16 unidades"
$ws.Cells.Item(74, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(74, 16).Value = 500
$ws.Cells.Item(74, 17).Value = 16
$ws.Cells.Item(74, 18).Value = "Hortaliza"
